# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (E2) and
# "Correspond Handback DateTime" (H2) timestamps on the per-locale
# report sheets, as produced by a fresh run of the handback-status
# report generator.

$wb = $excel.ActiveWorkbook

$ws_zhcn = $wb.Worksheets.Item("zh-cn")
$ws_zhcn.Range("E2").Value = "2016-03-18 10:50:01"
$ws_zhcn.Range("H2").Value = "2016-03-18 10:50:19"

$ws_dede = $wb.Worksheets.Item("de-de")
$ws_dede.Range("E2").Value = "2016-03-18 10:50:07"
$ws_dede.Range("H2").Value = "2016-03-18 10:50:24"
